# Update TC model prediction columns F:I (rows 2-41) with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 180.231689453125, 75.60057067871094, -485.054443359375, 26.0262279510498),
    @(3, 179.9908294677734, 71.97225952148438, -474.9625244140625, 24.97592163085938),
    @(4, 179.5555572509766, 68.16736602783203, -440.0684814453125, 49.60869979858398),
    @(5, 178.9087982177734, 64.22762298583984, -379.0177917480469, 127.9595565795898),
    @(6, 178.0307769775391, 60.20087814331055, -303.9894104003906, 185.3554382324219),
    @(7, 176.8992767333984, 56.13992309570312, -215.4222412109375, 211.9803466796875),
    @(8, 175.4897613525391, 52.10022735595703, -127.8701553344727, 230.8497924804688),
    @(9, 173.7764892578125, 48.13663864135742, -58.15444946289062, 244.5735626220703),
    @(10, 171.7340698242188, 44.29865264892578, -8.796058654785156, 249.2438659667969),
    @(11, 169.3401641845703, 40.62529373168945, 24.68533325195312, 252.4713439941406),
    @(12, 166.5794372558594, 37.14040756225586, 40.39848327636719, 251.6121826171875),
    @(13, 163.4480133056641, 33.8506965637207, 38.69892501831055, 239.6207580566406),
    @(14, 159.9584503173828, 30.74818229675293, 26.39996528625488, 220.3340301513672),
    @(15, 156.1432037353516, 27.81741142272949, 15.37427139282227, 203.8133850097656),
    @(16, 152.0557403564453, 25.04518508911133, 8.030043601989746, 192.6739807128906),
    @(17, 147.7683258056641, 22.42837905883789, 1.288897156715393, 181.418701171875),
    @(18, 143.3658599853516, 19.97669219970703, -4.337436676025391, 170.5476684570312),
    @(19, 138.9378967285156, 17.70913505554199, -8.467507362365723, 161.9538726806641),
    @(20, 134.5699920654297, 15.64644432067871, -12.68798065185547, 154.6984405517578),
    @(21, 130.3369140625, 13.80323028564453, -17.3156909942627, 147.717529296875),
    @(22, 126.2980651855469, 12.18283653259277, -21.78987884521484, 140.3986511230469),
    @(23, 122.4960327148438, 10.77618885040283, -25.7661190032959, 132.6969146728516),
    @(24, 118.9574279785156, 9.564249038696289, -29.17549133300781, 124.7720108032227),
    @(25, 115.6950302124023, 8.522297859191895, -32.02054977416992, 116.6630172729492),
    @(26, 112.7107238769531, 7.624358654022217, -34.40034484863281, 108.3360900878906),
    @(27, 109.9983444213867, 6.846392631530762, -36.51168441772461, 99.79259490966797),
    @(28, 107.5461730957031, 6.167929649353027, -38.58074569702148, 91.138427734375),
    @(29, 105.3391342163086, 5.572481155395508, -40.80999374389648, 82.60395050048828),
    @(30, 103.360107421875, 5.047057151794434, -43.36225891113281, 74.49432373046875),
    @(31, 101.5911865234375, 4.581432819366455, -46.34334182739258, 67.10404968261719),
    @(32, 100.0144348144531, 4.167401313781738, -49.76768112182617, 60.64458465576172),
    @(33, 98.61231231689453, 3.798205137252808, -53.54203414916992, 55.20652770996094),
    @(34, 97.36811065673828, 3.468092679977417, -57.49200439453125, 50.76336288452148),
    @(35, 96.26604461669922, 3.172082662582397, -61.41377639770508, 47.21026611328125),
    @(36, 95.29150390625, 2.905832052230835, -65.11982727050781, 44.4134635925293),
    @(37, 94.43093109130859, 2.665526151657104, -68.46622467041016, 42.24627304077148),
    @(38, 93.67191314697266, 2.447858572006226, -71.36507415771484, 40.60018157958984),
    @(39, 93.00320434570312, 2.249981164932251, -73.79103851318359, 39.37821960449219),
    @(40, 92.41458892822266, 2.069477319717407, -75.77970886230469, 38.49137878417969),
    @(41, 91.89684295654297, 1.904317021369934, -77.41483306884766, 37.86310195922852)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]   # F: H1_new_pred_KFold_mean
    $ws.Cells.Item($r, 7).Value = $row[2]   # G: H1_new_pred_KFold_std
    $ws.Cells.Item($r, 8).Value = $row[3]   # H: C2_new_pred_KFold_mean
    $ws.Cells.Item($r, 9).Value = $row[4]   # I: C2_new_pred_KFold_std
}
